$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '87.063.28'
$ws.Range('E2').Value = '  +10.63%  '
$ws.Range('D3').Value = '3.332.14'
$ws.Range('E3').Value = '  +7.92%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = '218.91'
$ws.Range('E5').Value = '  +8.93%  '
$ws.Range('D6').Value = '637.55'
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('D7').Value = '0.320'
$ws.Range('E7').Value = '  +27.27%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('D9').Value = '0.615'
$ws.Range('E9').Value = '  +8.48%  '
$ws.Range('D10').Value = '3.337.94'
$ws.Range('E10').Value = '  +8.29%  '
$ws.Range('D11').Value = '0.609'
$ws.Range('E11').Value = '  +10.80%  '
$ws.Range('D12').Value = '0.0000274'
$ws.Range('E12').Value = '  +17.49%  '
$ws.Range('E13').Value = '  +2.50%  '
$ws.Range('D14').Value = '3.951.13'
$ws.Range('E14').Value = '  +7.42%  '
$ws.Range('D15').Value = '34.26'
$ws.Range('E15').Value = '  +13.79%  '
$ws.Range('D16').Value = '5.41'
$ws.Range('E16').Value = '  +7.03%  '
$ws.Range('D17').Value = '86.629.33'
$ws.Range('E17').Value = '  +9.89%  '
$ws.Range('D18').Value = '3.339.84'
$ws.Range('E18').Value = '  +7.18%  '
$ws.Range('B19').Value = 'SuiNetwork'
$ws.Range('C19').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D19').Value = '3.24'
$ws.Range('E19').Value = '  +14.85%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '14.67'
$ws.Range('E20').Value = '  +6.88%  '
$ws.Range('D21').Value = '448.59'
$ws.Range('E21').Value = '  +6.43%  '
$ws.Range('D22').Value = '9.17'
$ws.Range('E22').Value = '  +4.89%  '
$ws.Range('D23').Value = '5.30'
$ws.Range('E23').Value = '  +5.95%  '
$ws.Range('D24').Value = '7.44'
$ws.Range('E24').Value = '  +9.75%  '
$ws.Range('D25').Value = '5.28'
$ws.Range('E25').Value = '  +18.00%  '
$ws.Range('D26').Value = '12.33'
$ws.Range('E26').Value = '  +19.43%  '
$ws.Range('D27').Value = '3.500.80'
$ws.Range('E27').Value = '  +6.45%  '
$ws.Range('D28').Value = '78.72'
$ws.Range('E28').Value = '  +6.01%  '
$ws.Range('D29').Value = '0.0000130'
$ws.Range('E29').Value = '  +13.23%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = '0.187'
$ws.Range('E31').Value = '  +57.43%  '
$ws.Range('D32').Value = '603.94'
$ws.Range('E32').Value = '  +13.50%  '
$ws.Range('D33').Value = '9.32'
$ws.Range('E33').Value = '  +8.23%  '
$ws.Range('D34').Value = '0.997'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').Value = '1.55'
$ws.Range('E35').Value = '  +10.66%  '
$ws.Range('E36').Value = '  +6.41%  '
$ws.Range('E37').Value = '  +6.91%  '
$ws.Range('D38').Value = '23.48'
$ws.Range('E38').Value = '  +6.49%  '
$ws.Range('D39').Value = '6.60'
$ws.Range('E39').Value = '  +24.07%  '
$ws.Range('D40').Value = '0.420'
$ws.Range('E40').Value = '  +8.87%  '
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = '21.36'
$ws.Range('E42').Value = '  +3.35%  '
$ws.Range('D43').Value = '3.10'
$ws.Range('E43').Value = '  +21.14%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Value = '2.06'
$ws.Range('E44').Value = '  +18.95%  '
$ws.Range('E45').Value = '  +0.17%  '
$ws.Range('D46').Value = '157.27'
$ws.Range('E46').Value = '  -2.93%  '
$ws.Range('D47').Value = '190.15'
$ws.Range('E47').Value = '  +4.46%  '
$ws.Range('E48').Value = '  +10.61%  '
$ws.Range('D49').Value = '45.72'
$ws.Range('E49').Value = '  +8.12%  '
$ws.Range('D50').Value = '0.787'
$ws.Range('E50').Value = '  +4.28%  '
$ws.Range('D51').Value = '0.657'
$ws.Range('E51').Value = '  +9.58%  '
